# Commit: "Rename example to correct layers"
#
# The example organizations import sheet used the human readable value
# "water board" in the "Layer" column. The application code expects
# machine-friendly, underscore separated layer identifiers, so the example
# data is corrected to use "water_board" instead of "water board".
#
# The workbook's styling was also refreshed: the header/data range got a
# white background fill (previously it had none), which is applied here too.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$usedRange = $ws.UsedRange
$rowCount = $usedRange.Rows.Count
$colCount = $usedRange.Columns.Count

# Find the "Layer" column by reading the header row, then fix every
# occurrence of "water board" -> "water_board" within that column.
$layerCol = 0
for ($c = 1; $c -le $colCount; $c++) {
    $header = $ws.Cells.Item(1, $c).Value2
    if ($header -eq "Layer") {
        $layerCol = $c
    }
}

if ($layerCol -eq 0) {
    $layerCol = 4
}

for ($r = 2; $r -le $rowCount; $r++) {
    $cell = $ws.Cells.Item($r, $layerCol)
    if ($cell.Value2 -eq "water board") {
        $cell.Value2 = "water_board"
    }
}

# Apply a plain white background fill across the whole used range
# (header row + data rows), matching the refreshed table styling.
$usedRange.Interior.Color = 16777215
